# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, newValue) updates for column F
$updates = @{
    "展览"     = @{ 2 = 4387; 3 = 2459; 6 = 53; 9 = 134; 10 = 155; 12 = 1641; 14 = 3491; 15 = 232 }
    "全部类型" = @{ 2 = 4387; 3 = 2459; 7 = 53; 11 = 134; 12 = 155; 16 = 1641; 18 = 3491; 19 = 232 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsForSheet = $updates[$sheetName]
    foreach ($row in $rowsForSheet.Keys) {
        $newValue = $rowsForSheet[$row]
        $ws.Cells.Item($row, 6).Value = $newValue
    }
}
